$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Estado de Cuenta" worker-period table (rows 16-36) is being refreshed:
#  - JOSE DE JESUS NAVARRO (doc 9692332) now lists first, periods 1703-1712
#    ascending, Salario Basico 781242.
#  - WILFREDO PEREZ HERAZO (doc 9162606) now lists second, periods 1801-1811
#    ascending, Salario Basico 828116 (Valor Mora 29509 for periods
#    1801-1808, 31249 for 1809-1811).
# ONESIMO ORTIZ MARTINEZ (row 37) is untouched.

$rows = @(
    @("9692332", "JOSE DE JESUS NAVARRO", "1703", 29509, 781242),
    @("9692332", "JOSE DE JESUS NAVARRO", "1704", 29509, 781242),
    @("9692332", "JOSE DE JESUS NAVARRO", "1705", 29509, 781242),
    @("9692332", "JOSE DE JESUS NAVARRO", "1706", 29509, 781242),
    @("9692332", "JOSE DE JESUS NAVARRO", "1707", 29509, 781242),
    @("9692332", "JOSE DE JESUS NAVARRO", "1708", 29509, 781242),
    @("9692332", "JOSE DE JESUS NAVARRO", "1709", 29509, 781242),
    @("9692332", "JOSE DE JESUS NAVARRO", "1710", 29509, 781242),
    @("9692332", "JOSE DE JESUS NAVARRO", "1711", 29509, 781242),
    @("9692332", "JOSE DE JESUS NAVARRO", "1712", 29509, 781242),
    @("9162606", "WILFREDO PEREZ HERAZO", "1801", 29509, 828116),
    @("9162606", "WILFREDO PEREZ HERAZO", "1802", 29509, 828116),
    @("9162606", "WILFREDO PEREZ HERAZO", "1803", 29509, 828116),
    @("9162606", "WILFREDO PEREZ HERAZO", "1804", 29509, 828116),
    @("9162606", "WILFREDO PEREZ HERAZO", "1805", 29509, 828116),
    @("9162606", "WILFREDO PEREZ HERAZO", "1806", 29509, 828116),
    @("9162606", "WILFREDO PEREZ HERAZO", "1807", 29509, 828116),
    @("9162606", "WILFREDO PEREZ HERAZO", "1808", 29509, 828116),
    @("9162606", "WILFREDO PEREZ HERAZO", "1809", 31249, 828116),
    @("9162606", "WILFREDO PEREZ HERAZO", "1810", 31249, 828116),
    @("9162606", "WILFREDO PEREZ HERAZO", "1811", 31249, 828116)
)

$r = 16
foreach ($row in $rows) {
    $ws.Cells.Item($r, 3).Value = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
    $ws.Cells.Item($r, 6).Value = $row[3]
    $ws.Cells.Item($r, 7).Value = $row[4]
    $r = $r + 1
}
